$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update patient record data (row 2)
$ws.Range("A2").Value = "Phil"
$ws.Range("B2").Value = "Salt"
$ws.Range("C2").Value = "phil@s.com"
$ws.Range("D2").Value = 3328293043
$ws.Range("E2").Value = "11/27/1989"

# Update hyperlink display text for the email cell (same target address,
# just a new friendly display text). Drop the existing link and recreate it
# through the range's own collection so the cell keeps its original look
# (re-creating via Hyperlinks.Add re-applies the default hyperlink style).
$c2 = $ws.Range("C2")
$c2.Hyperlinks.Delete()
$h = $c2.Hyperlinks.Item(1)
$h.Address = "mailto:Tim@d.com"
$h.TextToDisplay = "phil@s.com"

# Update the active selection to E2
$ws.Range("E2").Select()
